$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# list with the latest scraped figures. Numeric-looking Price values are
# written with a leading apostrophe so Excel stores them as text (matching
# the sheet's existing convention), exactly as a user typing them in
# would.

$ws.Range("D2").Value = "45.978.14"
$ws.Range("E2").Value = "  +3.33%  "

$ws.Range("D3").Value = "2.443.72"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'321.61"
$ws.Range("E5").Value = "  +2.66%  "

$ws.Range("D6").Value = "'104.36"
$ws.Range("E6").Value = "  +2.89%  "

$ws.Range("E7").Value = "  +1.04%  "

$ws.Range("D9").Value = "'0.536"
$ws.Range("E9").Value = "  +4.27%  "

$ws.Range("D10").Value = "'35.85"
$ws.Range("E10").Value = "  +1.86%  "

$ws.Range("D11").Value = "'0.0805"
$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("E12").Value = "  -2.07%  "

$ws.Range("D13").Value = "'18.21"
$ws.Range("E13").Value = "  -3.59%  "

$ws.Range("D14").Value = "'7.04"
$ws.Range("E14").Value = "  +1.61%  "

$ws.Range("D15").Value = "2.828.48"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").Value = "2.426.94"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").Value = "'0.839"
$ws.Range("E17").Value = "  +0.50%  "

$ws.Range("D18").Value = "45.801.28"
$ws.Range("E18").Value = "  +3.15%  "

$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").Value = "'6.41"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").Value = "0.0₃0932"
$ws.Range("E21").Value = "  +2.69%  "

$ws.Range("D22").Value = "'71.22"
$ws.Range("E22").Value = "  +3.46%  "

$ws.Range("E23").Value = "  +4.11%  "

$ws.Range("D24").Value = "'246.81"
$ws.Range("E24").Value = "  +2.58%  "

$ws.Range("E25").Value = "  +1.25%  "

$ws.Range("D27").Value = "'25.85"
$ws.Range("E27").Value = "  +2.66%  "

$ws.Range("D28").Value = "'2.20"
$ws.Range("E28").Value = "  -3.11%  "

$ws.Range("D29").Value = "'9.66"
$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").Value = "'33.57"
$ws.Range("E30").Value = "  +0.80%  "

$ws.Range("D31").Value = "'49.22"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("E32").Value = "  +2.59%  "

$ws.Range("D33").Value = "'20.01"
$ws.Range("E33").Value = "  +2.39%  "

$ws.Range("D34").Value = "'5.35"
$ws.Range("E34").Value = "  +3.26%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  -1.15%  "

$ws.Range("D37").Value = "'4.52"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("D39").Value = "'2.94"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").Value = "'126.15"
$ws.Range("E40").Value = "  +1.39%  "

$ws.Range("E41").Value = "  +3.89%  "

$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").Value = "'20.80"
$ws.Range("E43").Value = "  -3.32%  "

$ws.Range("E44").Value = "  +1.27%  "

$ws.Range("D45").Value = "1.957.23"
$ws.Range("E45").Value = "  +0.56%  "

$ws.Range("D46").Value = "'2.11"
$ws.Range("E46").Value = "  -3.05%  "

$ws.Range("E47").Value = "  +0.45%  "

$ws.Range("E48").Value = "  +10.87%  "

$ws.Range("D49").Value = "'9.09"
$ws.Range("E49").Value = "  -4.75%  "

$ws.Range("E50").Value = "  +7.51%  "

$ws.Range("D51").Value = "'77.27"
$ws.Range("E51").Value = "  +5.08%  "
